# Appends 9 new survey response rows (398-406) to the "Form Responses 1"
# sheet, mirroring the style of the last existing data row (397).
#
# Column map (row 1 headers):
#   A Timestamp (date/time serial, style 3)
#   B Kaupunki (city)
#   C Ika (age bracket)
#   D Sukupuoli (gender)
#   E Tyokokemus (years of experience)
#   F Tyosuhteen luonne (employment type)
#   G Tyoaika (work-time fraction, percentage style 4)
#   H Rooli (role)
#   I Etana vai paikallisesti (remote/on-site)
#   J Kuukausipalkka (monthly salary)
#   K Vuositulot / Vuosilaskutus (annual income/billing)
#   L Onko palkka kilpailukykyinen (salary competitive?)
#   M Tyopaikka (employer) -- unused by the new rows
#   N Vapaa sana (free text) -- only row 398 uses this

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet
$xlPasteFormats = -4122

$rows = @(
    @{ Row = 398; A = 44246.55214344907;  B = "Turku";    C = "31-35 v"; D = "Mies";
       E = 8;  F = "Työntekijä / palkollinen"; G = 1;   H = "Senior Software Engineer (Backend)";
       I = "Pääosin tai kokonaan etätyö"; J = 5700; K = 74100; L = "Kyllä";
       N = "Ennen koronaa oli osittainen etätyö, koronan jälkeen 100%" },

    @{ Row = 399; A = 44246.55232758102;  B = "EU";       C = "31-35 v"; D = "Mies";
       E = 3;  F = "Työntekijä / palkollinen"; G = 1;   H = "Ohjelmistokehittäjä";
       I = "Pääosin tai kokonaan toimistolla"; J = 3200; K = 40000; L = "Ei" },

    @{ Row = 400; A = 44246.56990268518;  B = "Jyväskylä"; C = "31-35 v"; D = "äiä";
       E = 6;  F = "Työntekijä / palkollinen"; G = 1;   H = "WordPress / Frontend-koodari";
       I = "Pääosin tai kokonaan etätyö"; J = 3000; K = "37 500"; L = "Kyllä" },

    @{ Row = 401; A = 44246.58209564815;  B = "Jyväskylä"; C = "21-25 v"; D = "Mies";
       E = 21; F = "Työntekijä / palkollinen"; G = 1;   H = "Arkkitehti";
       I = "Noin 50/50 hybridimalli"; J = 5500; K = 75000; L = "Kyllä" },

    @{ Row = 402; A = 44246.58394819444;  B = "PK-Seutu (Helsinki, Espoo, Vantaa)"; C = "31-35 v"; D = "Mies";
       E = 7;  F = "Työntekijä / palkollinen"; G = 1;   H = "Ohjelmistokehittäjä, backend";
       I = "Noin 50/50 hybridimalli"; J = 5470; K = 94000; L = "Kyllä" },

    @{ Row = 403; A = 44246.58876788194;  B = "PK-Seutu (Helsinki, Espoo, Vantaa)"; C = "31-35 v";
       E = 3;  F = "Työntekijä / palkollinen"; G = 1;   H = "Full stack developer";
       I = "Pääosin tai kokonaan etätyö"; J = 5300; K = 66250; L = "Kyllä" },

    @{ Row = 404; A = 44246.5905702662;   B = "Kuopio";   C = "31-35 v"; D = "Mies";
       E = 9;  F = "Työntekijä / palkollinen"; G = 0.8; H = "CTO";
       I = "Pääosin tai kokonaan etätyö"; J = 5200; K = 65000; L = "Kyllä" },

    @{ Row = 405; A = 44246.591665266205; B = "PK-Seutu (Helsinki, Espoo, Vantaa)"; C = "36-40 v"; D = "Mies";
       E = 14; F = "Työntekijä / palkollinen"; G = 1;   H = "Projektipäällikkö";
       I = "Noin 50/50 hybridimalli"; J = 6400; K = 102000; L = "Kyllä" },

    @{ Row = 406; A = 44246.60321177084;  B = "PK-Seutu (Helsinki, Espoo, Vantaa)"; C = "36-40 v"; D = "Mies";
       E = 15; F = "Työntekijä / palkollinen"; G = 1;   H = "Frontend & UX";
       I = "Pääosin tai kokonaan toimistolla"; J = 5000; K = "Optiot"; L = "Ei" }
)

foreach ($r in $rows) {
    $row = $r.Row

    # Clone number formatting from row 397 (A, G, and the general B:L cells)
    # so the new cells reuse the existing style indices (date / percentage /
    # general) instead of minting new ones in styles.xml. Only the columns
    # that will actually receive a value are formatted, so columns with no
    # data for this particular row (e.g. D on row 403) are left untouched
    # instead of becoming an empty-but-styled cell.
    $ws.Range("A397").Copy()
    $ws.Range("A$row").PasteSpecial($xlPasteFormats)
    $ws.Range("A$row").Value = $r.A

    $ws.Range("B397").Copy()
    $ws.Range("B$row").PasteSpecial($xlPasteFormats)
    $ws.Range("B$row").Value = $r.B

    if ($r.ContainsKey("C")) {
        $ws.Range("C397").Copy()
        $ws.Range("C$row").PasteSpecial($xlPasteFormats)
        $ws.Range("C$row").Value = $r.C
    }
    if ($r.ContainsKey("D")) {
        $ws.Range("D397").Copy()
        $ws.Range("D$row").PasteSpecial($xlPasteFormats)
        $ws.Range("D$row").Value = $r.D
    }

    $ws.Range("E397").Copy()
    $ws.Range("E$row").PasteSpecial($xlPasteFormats)
    $ws.Range("E$row").Value = $r.E

    $ws.Range("F397").Copy()
    $ws.Range("F$row").PasteSpecial($xlPasteFormats)
    $ws.Range("F$row").Value = $r.F

    $ws.Range("G397").Copy()
    $ws.Range("G$row").PasteSpecial($xlPasteFormats)
    $ws.Range("G$row").Value = $r.G

    $ws.Range("H397").Copy()
    $ws.Range("H$row").PasteSpecial($xlPasteFormats)
    $ws.Range("H$row").Value = $r.H

    $ws.Range("I397").Copy()
    $ws.Range("I$row").PasteSpecial($xlPasteFormats)
    $ws.Range("I$row").Value = $r.I

    $ws.Range("J397").Copy()
    $ws.Range("J$row").PasteSpecial($xlPasteFormats)
    $ws.Range("J$row").Value = $r.J

    $ws.Range("K397").Copy()
    $ws.Range("K$row").PasteSpecial($xlPasteFormats)
    $ws.Range("K$row").Value = $r.K

    $ws.Range("L397").Copy()
    $ws.Range("L$row").PasteSpecial($xlPasteFormats)
    $ws.Range("L$row").Value = $r.L

    if ($r.ContainsKey("N")) {
        $ws.Range("N23").Copy()
        $ws.Range("N$row").PasteSpecial($xlPasteFormats)
        $ws.Range("N$row").Value = $r.N
    }
}
